$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, F, J, K across rows 2-25
# (loading_percent results for the 380 kV case)
$newValues = @(
    ,@(15.38225564643492, 7.466367187508107, 4.458777462262151, 11.36521475653977, 60.79824859264363, 10.83153443283537, 15.0308099310607)
    ,@(15.31761715088765, 7.461854550307004, 4.511138247076429, 11.41148789700232, 59.82695907145634, 10.82488498585519, 15.02419668091847)
    ,@(15.28389578465265, 7.463283789944406, 4.545934970892977, 11.44350797748967, 59.22646136411002, 10.82240464489274, 15.02566999335046)
    ,@(15.2716666762217, 7.464919303484549, 4.560773955365147, 11.45746059880429, 58.98092973084075, 10.82179724639128, 15.02766123609725)
    ,@(15.26972770708928, 7.465254364239289, 4.56327757370855, 11.459831937166, 58.94011588263041, 10.82172075898999, 15.0280758389538)
    ,@(15.28372471993663, 7.463301588364042, 4.546132434777301, 11.44369249070242, 59.22315309243624, 10.82239481969377, 15.02569121839029)
    ,@(15.35873794015031, 7.463938159160707, 4.476278061671829, 11.38041901260136, 60.46433016239355, 10.82890940980146, 15.02738111146379)
    ,@(15.55252708302363, 7.498571413083943, 4.360601037299372, 11.28510406073274, 62.85640637024959, 10.8543758215049, 15.07457440308378)
    ,@(15.72231997927095, 7.544339504150178, 4.289035790225433, 11.23278799735229, 64.57641039704164, 10.88078393743061, 15.13587063765806)
    ,@(15.80524617340496, 7.569532995982827, 4.25949170523734, 11.21286864452582, 65.34834109168392, 10.89445734332768, 15.16948104229567)
    ,@(15.83744016758567, 7.579696122334028, 4.248745649393559, 11.20588613188389, 65.6389477808508, 10.89987247404824, 15.18302496355924)
    ,@(15.83047184032013, 7.577479719631993, 4.251040222347366, 11.20736497179286, 65.57643936396339, 10.898695699098, 15.18007185548379)
    ,@(15.80787908262969, 7.570356678564753, 4.258598719172666, 11.21228294099462, 65.37228477515291, 10.89489810220806, 15.17057899245966)
    ,@(15.79414265725843, 7.566074525461001, 4.263286292624326, 11.21536840989533, 65.24700626467302, 10.8926028235123, 15.16487042503289)
    ,@(15.71701279791642, 7.54278057770525, 4.291027892412172, 11.23416806987976, 64.52573454917517, 10.8799236381224, 15.13378877651082)
    ,@(15.67113462143466, 7.529606861266271, 4.308823789676511, 11.24669662688747, 64.08042457533745, 10.87256963262127, 15.11618327811742)
    ,@(15.6452836362636, 7.522441789788897, 4.319342425389312, 11.25426773070421, 63.82331621925609, 10.86849616835236, 15.10659657147558)
    ,@(15.63662391933686, 7.520086754987688, 4.322952185413928, 11.25689378317225, 63.73610190117308, 10.86714385974187, 15.10344353156671)
    ,@(15.67596304473776, 7.530966613262517, 4.306900044634059, 11.24532514375555, 64.12793122869311, 10.87333630563099, 15.11800162306211)
    ,@(15.81449385483713, 7.572432038595073, 4.256366549242196, 11.21082318189084, 65.43229776446232, 10.89600712027866, 15.1733451828598)
    ,@(15.9096327503991, 7.603158946629422, 4.225918429153194, 11.19154234931833, 66.27474392496411, 10.91220623311811, 15.21426972055903)
    ,@(15.85844340276877, 7.586429939292691, 4.241930321831142, 11.20153302806605, 65.82609501794747, 10.90343448904513, 15.19199518623187)
    ,@(15.67377847738079, 7.530350596022761, 4.3077688748034, 11.24594404440501, 64.10645686246903, 10.87298921157516, 15.11717788260899)
    ,@(15.49520254664394, 7.485624713897806, 4.389569891709205, 11.3077910515199, 62.21502169856129, 10.8543758215049, 15.07457440308378)
)

$cols = @(2, 3, 4, 5, 6, 10, 11)
$startRow = 2

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $colNum = $cols[$j]
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j]
    }
}

Write-Output "Updated $($newValues.Count) rows across columns B,C,D,E,F,J,K"
